$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.180.83'
$ws.Range("E2").Value = '  -3.19%  '
$ws.Range("D3").Value = '1.714.86'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.66'
$ws.Range("E5").Value = '  -5.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4788'
$ws.Range("E7").Value = '  +6.30%  '
$ws.Range("E8").Value = '  -3.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.27'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07281'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.047'
$ws.Range("E11").Value = '  -5.67%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.89'
$ws.Range("E13").Value = '  -5.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.868'
$ws.Range("E14").Value = '  -3.32%  '
$ws.Range("D15").Value = '1.713.09'
$ws.Range("E15").Value = '  -3.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.880'
$ws.Range("E16").Value = '  -5.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.02'
$ws.Range("E17").Value = '  -5.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001041'
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06357'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.50'
$ws.Range("E21").Value = '  -3.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.627'
$ws.Range("E22").Value = '  -3.14%  '
$ws.Range("D23").Value = '27.207.84'
$ws.Range("E23").Value = '  -3.13%  '
$ws.Range("E24").Value = '  -4.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.088'
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.99'
$ws.Range("E26").Value = '  -6.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.70'
$ws.Range("E27").Value = '  -3.60%  '
$ws.Range("D28").Value = '1.909.96'
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("E29").Value = '  -3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.05'
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.020'
$ws.Range("E31").Value = '  -8.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09235'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.580'
$ws.Range("E33").Value = '  -3.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.337'
$ws.Range("E34").Value = '  -6.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02204'
$ws.Range("E35").Value = '  -4.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05924'
$ws.Range("E36").Value = '  -4.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.10'
$ws.Range("E37").Value = '  -6.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2007'
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.420'
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.758'
$ws.Range("E40").Value = '  -5.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5948'
$ws.Range("E42").Value = '  -6.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.098'
$ws.Range("E43").Value = '  -7.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.489'
$ws.Range("E44").Value = '  -5.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.62'
$ws.Range("E45").Value = '  -5.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.583'
$ws.Range("E46").Value = '  -4.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5628'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.76'
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.842'
$ws.Range("E49").Value = '  -6.32%  '
$ws.Range("E50").Value = '  -3.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.086'
$ws.Range("E51").Value = '  -5.18%  '
